$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "616.51", "1.00", "0.0000248") are stored as
# text, matching the workbook's existing inlineStr cells instead of being
# auto-converted to numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '90.812.34'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '3.080.38'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '243.69'
$ws.Range('E5').Value = '  +2.70%  '
$ws.Range('D6').Value = '616.51'
$ws.Range('E6').Value = '  -2.10%  '
$ws.Range('E7').Value = '  +7.77%  '
$ws.Range('D8').Value = '0.365'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = '3.075.41'
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('D11').Value = '0.740'
$ws.Range('E11').Value = '  +3.60%  '
$ws.Range('E12').Value = '  +2.96%  '
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').Value = '34.95'
$ws.Range('E14').Value = '  -4.36%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '91.039.43'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('B16').Value = 'Toncoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D16').Value = '5.42'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Value = '3.661.35'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '3.096.93'
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('E19').Value = '  -2.59%  '
$ws.Range('D20').Value = '14.43'
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('D21').Value = '0.0000211'
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('D22').Value = '5.72'
$ws.Range('E22').Value = '  +3.20%  '
$ws.Range('D23').Value = '440.24'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').Value = '9.02'
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('D25').Value = '90.93'
$ws.Range('E25').Value = '  +3.57%  '
$ws.Range('D26').Value = '5.59'
$ws.Range('E26').Value = '  -5.23%  '
$ws.Range('D27').Value = '11.75'
$ws.Range('E27').Value = '  -5.70%  '
$ws.Range('D28').Value = '3.260.27'
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.248'
$ws.Range('E30').Value = '  +27.64%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').Value = '0.183'
$ws.Range('E31').Value = '  +15.06%  '
$ws.Range('D32').Value = '9.13'
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('E33').Value = '  +13.69%  '
$ws.Range('E34').Value = '  +14.30%  '
$ws.Range('D35').Value = '0.110'
$ws.Range('E35').Value = '  +30.60%  '
$ws.Range('D36').Value = '7.69'
$ws.Range('E36').Value = '  +8.41%  '
$ws.Range('D37').Value = '26.31'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = '4.19'
$ws.Range('E38').Value = '  +29.17%  '
$ws.Range('D39').Value = '1.91'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = '490.70'
$ws.Range('E40').Value = '  -3.33%  '
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').Value = '  -4.79%  '
$ws.Range('D42').Value = '1.28'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = '0.415'
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('D44').Value = '22.12'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = '154.01'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '1.88'
$ws.Range('E47').Value = '  -0.76%  '
$ws.Range('D48').Value = '0.682'
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('D49').Value = '4.44'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').Value = '44.06'
$ws.Range('E51').Value = '  -2.52%  '
